$d = $word.ActiveDocument

# Update the delivery date
$d.Content.Find.Execute("07/07", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "01/12", 2)

# Change the group-size sentence: "grupos de 3 pessoas" -> "grupos"
$d.Content.Find.Execute("grupos de 3 pessoas", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "grupos", 2)

# Update the reference from Trabalho 1 to Trabalho 2
$d.Content.Find.Execute("mesmos definidos no Trabalho 1.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "mesmos definidos no Trabalho 2.", 2)
